$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.772.22'
$ws.Range("E2").Value = '  +1.64%  '
$ws.Range("D3").Value = '2.491.60'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.13%  '
$ws.Range("E7").Value = '  +0.26%  '
$ws.Range("E8").Value = '  +3.42%  '
$ws.Range("D9").Value = '2.499.90'
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("E11").Value = '  -2.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.19%  '
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '2.933.13'
$ws.Range("E14").Value = '  +1.65%  '
$ws.Range("D15").Value = '58.599.99'
$ws.Range("E15").Value = '  +1.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.69%  '
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = '2.490.87'
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("E22").Value = '  +4.87%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.03%  '
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.993'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.23%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.159'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.84'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.65%  '
$ws.Range("E30").Value = '  +3.54%  '
$ws.Range("E31").Value = '  +3.94%  '
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("E33").Value = '  +1.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.70%  '
$ws.Range("E37").Value = '  -2.93%  '
$ws.Range("E38").Value = '  +1.31%  '
$ws.Range("E39").Value = '  +3.77%  '
$ws.Range("B40").Value = 'SuiNetwork'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.821'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.65%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.32'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.46'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.89%  '
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.09'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '273.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '131.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.591'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.31%  '
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0508'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.20%  '
$ws.Range("E49").Value = '  +3.45%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.52'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.58%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.58%  '
